# Elite Support Datasheet - small copy-edit on the "Support Hours" table
# (slide sldId=261 / cId=1050037809, graphicFrame id=25, name "Table 6"):
#   - "...English and Japanese " -> "...English and Japanese. " (nbsp kept)
#   - "...only in Japan"         -> "...only in Japan."
#
# Locate the shape robustly (by Id + HasTable) instead of hard-coding a
# slide index, since several slides contain a shape literally named
# "Table 6".

$p = $ppt.ActivePresentation

$tableShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $cand = $slide.Shapes.Item($shi)
        if ($cand.Id -eq 25 -and $cand.HasTable) {
            $tableShape = $cand
        }
    }
}

$tbl = $tableShape.Table
$cell = $tbl.Cell(3, 1)
$tr = $cell.Shape.TextFrame.TextRange

$nbsp = [char]0x00A0

# Rebuild the cell's full text with the two wording fixes applied, keeping
# the existing paragraph breaks and the trailing non-breaking space intact.
$newText = "Language support is only available in English and Japanese." + $nbsp `
    + "`r" `
    + "`r" `
    + $nbsp + "1 P2, P3, P4 cases are limited to business hours only in Japan."

$tr.Text = $newText
